# Duplicate column A (the translation-key labels) into column B for rows 1-28,
# so each row now has the same text in both A and B (supports the v2 update
# translation patch which needs a second, parallel column of the same keys).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 28; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Text
}
